# Add season-record columns (Wins / Losses / Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the header style (bold, bordered, centered)
# from an existing header cell so the new header cells share style index 1,
# then set the header text.
$ws.Range("A1").Copy($ws.Range("AC1"))
$ws.Range("A1").Copy($ws.Range("AD1"))
$ws.Range("A1").Copy($ws.Range("AE1"))

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# --- Data rows (2-50): every player on this roster shares the team's
# 1992 season record: 70 wins, 92 losses, 0 ties.
$ws.Range("AC2:AC50").Value = 70
$ws.Range("AD2:AD50").Value = 92
$ws.Range("AE2:AE50").Value = 0

Write-Output "season record columns added"
